$p = $ppt.ActivePresentation

# --- Slide 1: remove the review comment and clear all shapes/content ---
$s1 = $p.Slides.Item(1)

while ($s1.Comments.Count -gt 0) {
    $s1.Comments.Item(1).Delete()
}

while ($s1.Shapes.Count -gt 0) {
    $s1.Shapes.Item(1).Delete()
}

# --- Slide 2: clear all shapes/content and blank out its name ---
$s2 = $p.Slides.Item(2)

while ($s2.Shapes.Count -gt 0) {
    $s2.Shapes.Item(1).Delete()
}

$s2.Name = ""

# --- Slide 3: remove entirely ---
$p.Slides.Item(3).Delete()
